$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add English / Traditional-Chinese translation columns (C = en, D = tw) ---
# Row 4 - Chang'an (B4 already holds the CN text "长安")
$ws.Range("C4").Value = "Chang'an"
$ws.Range("D4").Value = "長安"

# Row 5 - Aolai country
$ws.Range("B5").Value = "傲来国"
$ws.Range("C5").Value = "Aolai country"
$ws.Range("D5").Value = "傲來國"

# Row 6 - Spiritual mountain
$ws.Range("B6").Value = "灵山"
$ws.Range("C6").Value = "Spiritual mountain"
$ws.Range("D6").Value = "靈山"

# Row 7 - Tiangong (note trailing newline kept from source data)
$ws.Range("B7").Value = "天宫"
$ws.Range("C7").Value = "Tiangong`n"
$ws.Range("D7").Value = "天宮"

# Row 8 - underworld (tw column re-uses the cn text, matching source data)
$ws.Range("B8").Value = "地府"
$ws.Range("C8").Value = "underworld"
$ws.Range("D8").Value = "地府"

# Row 9 - Flame Mountain (tw column re-uses the cn text, matching source data)
$ws.Range("B9").Value = "火焰山"
$ws.Range("C9").Value = "Flame Mountain`n"
$ws.Range("D9").Value = "火焰山"

# --- Row heights ---
$ws.Range("A1:E1").RowHeight = 16.5
$ws.Range("A2:E2").RowHeight = 16.5
$ws.Range("A3:E3").RowHeight = 16.5
$ws.Range("A7:D7").RowHeight = 27
$ws.Range("A9:D9").RowHeight = 40.5

# --- Column widths (C holds the new "en" text, E already had a custom width) ---
$ws.Columns.Item(3).ColumnWidth = 20.375
$ws.Columns.Item(5).ColumnWidth = 17.6166666666667

# --- Wrap the long translations in C7 / C9 ---
$ws.Range("C7").WrapText = $true
$ws.Range("C7").VerticalAlignment = -4108
$ws.Range("C9").WrapText = $true
$ws.Range("C9").VerticalAlignment = -4108

# --- Selection / window state left by the author on save ---
$ws.Range("C12").Select() | Out-Null
$excel.ActiveWindow.Width = 27945
$excel.ActiveWindow.Height = 12255
